$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 - this shifts existing rows 26..72 down to 27..73
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly record
$ws.Cells.Item(26, 1).Value = 4
$ws.Cells.Item(26, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(26, 3).Value = "Los Lagos"
$ws.Cells.Item(26, 4).Value = 44979
$ws.Cells.Item(26, 5).Value = 10
$ws.Cells.Item(26, 6).Value = 100112043
$ws.Cells.Item(26, 7).Value = "Pepino dulce"
$ws.Cells.Item(26, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 25
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = 20000
$ws.Cells.Item(26, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(26, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(26, 16).Value = 1111
$ws.Cells.Item(26, 17).Value = 18
$ws.Cells.Item(26, 18).Value = "Hortaliza"
